$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 21003.5
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 14285750
$ws.Range("I8").Value = 14285750
$ws.Range("K8").Value = 42857250
$ws.Range("M8").Value = -42857111

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H14").Value = 21003.5
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2748.7273
$ws.Range("I19").Value = 1256.1666
$ws.Range("K19").Value = 1256.1666
$ws.Range("M19").Value = -1081.1666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 18753022
$ws.Range("I32").Value = 25000922
$ws.Range("K32").Value = 25000922
$ws.Range("M32").Value = -25000596

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H57").Value = 132524
$ws.Range("J57").Value = 132524
$ws.Range("L57").Value = 397572
$ws.Range("N57").Value = -398570

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 15072.417
$ws.Range("I69").Value = 6898.8
$ws.Range("J69").Value = 20910.715
$ws.Range("K69").Value = 20696.4
$ws.Range("L69").Value = 62732.145
$ws.Range("M69").Value = -19822.4
$ws.Range("N69").Value = -64480.145

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 15072.417
$ws.Range("I72").Value = 6898.8
$ws.Range("J72").Value = 20910.715
$ws.Range("K72").Value = 62089.2
$ws.Range("L72").Value = 188196.435
$ws.Range("M72").Value = -57721.2
$ws.Range("N72").Value = -196932.435

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 1516
$ws.Range("I76").Value = 1516
$ws.Range("K76").Value = 1516
$ws.Range("M76").Value = -1201

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 1516
$ws.Range("I79").Value = 1516
$ws.Range("K79").Value = 1516
$ws.Range("M79").Value = -424

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4964.6665
$ws.Range("I86").Value = 4388.5
$ws.Range("J86").Value = 6117
$ws.Range("K86").Value = 4388.5
$ws.Range("L86").Value = 6117
$ws.Range("M86").Value = -3265.5
$ws.Range("N86").Value = -8363

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 4964.6665
$ws.Range("I89").Value = 4388.5
$ws.Range("J89").Value = 6117
$ws.Range("K89").Value = 21942.5
$ws.Range("L89").Value = 30585
$ws.Range("M89").Value = -16326.5
$ws.Range("N89").Value = -41817

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1960.8948
$ws.Range("I100").Value = 1884.4286
$ws.Range("K100").Value = 1884.4286
$ws.Range("M100").Value = -1343.4286

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 9826.691999999999
$ws.Range("I113").Value = 10294.637
$ws.Range("K113").Value = 10294.637
$ws.Range("M113").Value = -7040.637000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 6754.2
$ws.Range("J132").Value = 9072.357
$ws.Range("L132").Value = 27217.071
$ws.Range("N132").Value = -32277.071

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2597.4524
$ws.Range("I138").Value = 1561
$ws.Range("K138").Value = 4683
$ws.Range("M138").Value = 457

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1100546.9
$ws.Range("I2").Value = 1362136.5
$ws.Range("J2").Value = 1870.8
$ws.Range("K2").Value = 1362136.5
$ws.Range("L2").Value = 1870.8
$ws.Range("M2").Value = -1362023.5
$ws.Range("N2").Value = -2096.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9258.134
$ws.Range("J32").Value = 24126.285
$ws.Range("L32").Value = 24126.285
$ws.Range("N32").Value = -24700.285

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H70").Value = 59999
$ws.Range("J70").Value = 59999
$ws.Range("L70").Value = 59999
$ws.Range("N70").Value = -60539

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H73").Value = 59999
$ws.Range("J73").Value = 59999
$ws.Range("L73").Value = 59999
$ws.Range("N73").Value = -61871

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1224.2727
$ws.Range("I74").Value = 900.36664
$ws.Range("J74").Value = 4463.3335
$ws.Range("K74").Value = 900.36664
$ws.Range("L74").Value = 4463.3335
$ws.Range("M74").Value = -26.36663999999996
$ws.Range("N74").Value = -6211.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1224.2727
$ws.Range("I77").Value = 900.36664
$ws.Range("J77").Value = 4463.3335
$ws.Range("K77").Value = 4501.8332
$ws.Range("L77").Value = 22316.6675
$ws.Range("M77").Value = -133.8332
$ws.Range("N77").Value = -31052.6675

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H96").Value = 30000
$ws.Range("J96").Value = 30000
$ws.Range("L96").Value = 30000
$ws.Range("N96").Value = -35492

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H101").Value = 48998
$ws.Range("J101").Value = 48998
$ws.Range("L101").Value = 48998
$ws.Range("N101").Value = -55488

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1100546.9
$ws.Range("I116").Value = 1362136.5
$ws.Range("J116").Value = 1870.8
$ws.Range("K116").Value = 1362136.5
$ws.Range("L116").Value = 1870.8
$ws.Range("M116").Value = -1359842.5
$ws.Range("N116").Value = -6458.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1811.6
$ws.Range("I122").Value = 1804.7727
$ws.Range("K122").Value = 5414.3181
$ws.Range("M122").Value = -2964.3181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1100546.9
$ws.Range("I3").Value = 1362136.5
$ws.Range("J3").Value = 1870.8
$ws.Range("K3").Value = 1362136.5
$ws.Range("L3").Value = 1870.8
$ws.Range("M3").Value = -1362022.5
$ws.Range("N3").Value = -2098.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5452.125
$ws.Range("I134").Value = 4689.4
$ws.Range("K134").Value = 14068.2
$ws.Range("M134").Value = -11533.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 10008
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2211.8667
$ws.Range("I31").Value = 1928.3043
$ws.Range("K31").Value = 1928.3043
$ws.Range("M31").Value = -1633.3043

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2211.8667
$ws.Range("I34").Value = 1928.3043
$ws.Range("K34").Value = 1928.3043
$ws.Range("M34").Value = -1726.3043

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 11021.549
$ws.Range("J99").Value = 12412.883
$ws.Range("L99").Value = 12412.883
$ws.Range("N99").Value = -15408.883

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1048.0667
$ws.Range("I105").Value = 1152.3334
$ws.Range("J105").Value = 804.7778
$ws.Range("K105").Value = 1152.3334
$ws.Range("L105").Value = 804.7778
$ws.Range("M105").Value = 594.6666
$ws.Range("N105").Value = -4298.7778

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 11021.549
$ws.Range("J126").Value = 12412.883
$ws.Range("L126").Value = 37238.649
$ws.Range("N126").Value = -42178.649

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 14545.091
$ws.Range("I132").Value = 15599.6
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 46798.8
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -44268.8
$ws.Range("N132").Value = -17060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 556.3333
$ws.Range("I44").Value = 269
$ws.Range("J44").Value = 700
$ws.Range("K44").Value = 807
$ws.Range("L44").Value = 2100
$ws.Range("M44").Value = -409
$ws.Range("N44").Value = -2896

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 16668332
$ws.Range("I3").Value = 10000000
$ws.Range("J3").Value = 20002498
$ws.Range("K3").Value = 10000000
$ws.Range("L3").Value = 20002498
$ws.Range("M3").Value = -9999884
$ws.Range("N3").Value = -20002730

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2276.5881
$ws.Range("I16").Value = 2276.5881
$ws.Range("K16").Value = 2276.5881
$ws.Range("M16").Value = -2106.5881

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 577.5
$ws.Range("I82").Value = 554.2967
$ws.Range("J82").Value = 999.8
$ws.Range("K82").Value = 554.2967
$ws.Range("L82").Value = 999.8
$ws.Range("M82").Value = -193.2967
$ws.Range("N82").Value = -1721.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 577.5
$ws.Range("I85").Value = 554.2967
$ws.Range("J85").Value = 999.8
$ws.Range("K85").Value = 554.2967
$ws.Range("L85").Value = 999.8
$ws.Range("M85").Value = 693.7033
$ws.Range("N85").Value = -3495.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2441.0715
$ws.Range("I93").Value = 2423.25
$ws.Range("J93").Value = 2548
$ws.Range("K93").Value = 2423.25
$ws.Range("L93").Value = 2548
$ws.Range("M93").Value = -1175.25
$ws.Range("N93").Value = -5044

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 8795.75
$ws.Range("J41").Value = 8729.833000000001
$ws.Range("L41").Value = 8729.833000000001
$ws.Range("N41").Value = -9509.833000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8561.625
$ws.Range("I62").Value = 8415.5
$ws.Range("K62").Value = 8415.5
$ws.Range("M62").Value = -7791.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 8561.625
$ws.Range("I65").Value = 8415.5
$ws.Range("K65").Value = 42077.5
$ws.Range("M65").Value = -38957.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 6677.8335
$ws.Range("I122").Value = 4469.64
$ws.Range("K122").Value = 13408.92
$ws.Range("M122").Value = -10958.92

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1537.9667
$ws.Range("I132").Value = 1221.4286
$ws.Range("J132").Value = 2276.5557
$ws.Range("K132").Value = 3664.2858
$ws.Range("L132").Value = 6829.6671
$ws.Range("M132").Value = -1134.2858
$ws.Range("N132").Value = -11889.6671
